$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 6000
$ws.Range("J32").Value = 6000
$ws.Range("L32").Value = 6000
$ws.Range("N32").Value = -6652
# Row 62
$ws.Range("H62").Value = 1099.6666
$ws.Range("I62").Value = 1099.6666
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 1099.6666
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -475.6666
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 1099.6666
$ws.Range("I65").Value = 1099.6666
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 5498.333000000001
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -2378.333000000001
$ws.Range("N65").ClearContents()
# Row 95
$ws.Range("H95").Value = 33116
$ws.Range("J95").Value = 33116
$ws.Range("L95").Value = 33116
$ws.Range("N95").Value = -38608
# Row 100
$ws.Range("H100").Value = 2531.6
$ws.Range("I100").Value = 3081.25
$ws.Range("J100").Value = 333
$ws.Range("K100").Value = 3081.25
$ws.Range("L100").Value = 333
$ws.Range("M100").Value = -2540.25
$ws.Range("N100").Value = -1415
# Row 101
$ws.Range("H101").Value = 200
$ws.Range("I101").Value = 200
$ws.Range("K101").Value = 600
$ws.Range("M101").Value = 1022
# Row 113
$ws.Range("H113").Value = 35961760
$ws.Range("I113").Value = 13892234
$ws.Range("K113").Value = 13892234
$ws.Range("M113").Value = -13888980
# Row 116
$ws.Range("H116").Value = 31262124
$ws.Range("I116").Value = 125005000
$ws.Range("K116").Value = 125005000
$ws.Range("M116").Value = -125001558
# Row 132
$ws.Range("H132").Value = 1814.1818
$ws.Range("I132").Value = 1346
$ws.Range("K132").Value = 4038
$ws.Range("M132").Value = -1508

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3577014.8
$ws.Range("J32").Value = 44999.5
$ws.Range("L32").Value = 44999.5
$ws.Range("N32").Value = -45573.5
# Row 45
$ws.Range("H45").Value = 3696.4167
$ws.Range("I45").Value = 2936.6667
$ws.Range("J45").Value = 3949.6667
$ws.Range("K45").Value = 2936.6667
$ws.Range("L45").Value = 3949.6667
$ws.Range("M45").Value = -2559.6667
$ws.Range("N45").Value = -4703.6667
# Row 86
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
# Row 89
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
# Row 138
$ws.Range("H138").Value = 99580
$ws.Range("J138").Value = 99580
$ws.Range("L138").Value = 99580
$ws.Range("N138").Value = -109860

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 18
$ws.Range("H18").Value = 1000
$ws.Range("I18").Value = 1000
$ws.Range("K18").Value = 1000
$ws.Range("M18").Value = -471
# Row 94
$ws.Range("H94").Value = 3782.6667
$ws.Range("J94").Value = 10002.5
$ws.Range("L94").Value = 10002.5
$ws.Range("N94").Value = -10904.5
# Row 104
$ws.Range("H104").Value = 67000
$ws.Range("J104").Value = 67000
$ws.Range("L104").Value = 67000
$ws.Range("N104").Value = -73988

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 5818.3228
$ws.Range("I16").Value = 706.5
$ws.Range("K16").Value = 706.5
$ws.Range("M16").Value = -419.5
# Row 31
$ws.Range("H31").Value = 6560.2085
$ws.Range("I31").Value = 1895.4546
$ws.Range("K31").Value = 1895.4546
$ws.Range("M31").Value = -1600.4546
# Row 34
$ws.Range("H34").Value = 6560.2085
$ws.Range("I34").Value = 1895.4546
$ws.Range("K34").Value = 1895.4546
$ws.Range("M34").Value = -1693.4546
# Row 113
$ws.Range("H113").Value = 5818.3228
$ws.Range("I113").Value = 706.5
$ws.Range("K113").Value = 706.5
$ws.Range("M113").Value = 1463.5
# Row 122
$ws.Range("H122").Value = 2812.8
$ws.Range("I122").Value = 1880.6666
$ws.Range("K122").Value = 5641.9998
$ws.Range("M122").Value = -3191.9998

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1867.9445
$ws.Range("J5").Value = 3834
$ws.Range("L5").Value = 11502
$ws.Range("N5").Value = -11726
# Row 36
$ws.Range("H36").Value = 500000100
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
# Row 113
$ws.Range("H113").Value = 776.1
$ws.Range("J113").Value = 816.7143
$ws.Range("L113").Value = 2450.1429
$ws.Range("N113").Value = -6790.1429
# Row 135
$ws.Range("H135").Value = 1867.9445
$ws.Range("J135").Value = 3834
$ws.Range("L135").Value = 34506
$ws.Range("N135").Value = -39576

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 54
$ws.Range("H54").Value = 32826
$ws.Range("J54").Value = 32826
$ws.Range("L54").Value = 32826
$ws.Range("N54").Value = -33606
# Row 95
$ws.Range("H95").Value = 24863.857
$ws.Range("J95").Value = 24863.857
$ws.Range("L95").Value = 24863.857
$ws.Range("N95").Value = -30355.857
# Row 98
$ws.Range("H98").Value = 43994.6
$ws.Range("J98").Value = 43994.6
$ws.Range("L98").Value = 43994.6
$ws.Range("N98").Value = -49984.6
# Row 100
$ws.Range("H100").Value = 51984.5
$ws.Range("J100").Value = 51984.5
$ws.Range("L100").Value = 51984.5
$ws.Range("N100").Value = -54148.5
# Row 101
$ws.Range("H101").Value = 53969
$ws.Range("J101").Value = 53969
$ws.Range("L101").Value = 53969
$ws.Range("N101").Value = -60459
# Row 122
$ws.Range("H122").Value = 7265806.5
$ws.Range("I122").Value = 7265806.5
$ws.Range("K122").Value = 21797419.5
$ws.Range("M122").Value = -21794969.5

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 9
$ws.Range("H9").Value = 147.5
$ws.Range("I9").Value = 146.66667
$ws.Range("K9").Value = 146.66667
$ws.Range("M9").Value = 77.33332999999999
# Row 122
$ws.Range("H122").Value = 4124.16
$ws.Range("I122").Value = 2682.8235
$ws.Range("K122").Value = 8048.470499999999
$ws.Range("M122").Value = -5598.470499999999

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
# Row 86
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
# Row 87
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
# Row 88
$ws.Range("H88").Value = 64689
$ws.Range("J88").Value = 64689
$ws.Range("L88").Value = 64689
$ws.Range("N88").Value = -65501
# Row 89
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
# Row 90
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
# Row 91
$ws.Range("H91").Value = 64689
$ws.Range("J91").Value = 64689
$ws.Range("L91").Value = 64689
$ws.Range("N91").Value = -67497
# Row 122
$ws.Range("H122").Value = 164475.72
$ws.Range("I122").Value = 225700.44
$ws.Range("K122").Value = 677101.3200000001
$ws.Range("M122").Value = -674651.3200000001
# Row 126
$ws.Range("H126").Value = 965.8333
$ws.Range("I126").Value = 995
$ws.Range("J126").Value = 960
$ws.Range("K126").Value = 2985
$ws.Range("L126").Value = 2880
$ws.Range("M126").Value = -515
$ws.Range("N126").Value = -7820
# Row 136
$ws.Range("H136").Value = 338327.44
$ws.Range("I136").Value = 1801.3334
$ws.Range("J136").Value = 562678.2
$ws.Range("K136").Value = 5404.0002
$ws.Range("L136").Value = 1688034.6
$ws.Range("M136").Value = -2854.0002
$ws.Range("N136").Value = -1693134.6
